$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column C. This shifts the existing column C
# (and its formatting/styles) two columns to the right, to column E.
$ws.Range("C1:D1").EntireColumn.Insert()

# Give the two newly inserted columns (C, D) the same custom width as column E
# (8.0 "points-ish" Excel column-width units -> serializes to width="8.0").
$ws.Columns("C:D").ColumnWidth = 7.1666666666666666

# Header row: B1 becomes the newest date header, C1 the next newest (both new),
# D1 keeps the date that used to be in B1, and E1 keeps the date that used to
# be in C1 (already shifted there automatically by the column insert above).
$ws.Range("B1").Value = "Jun_17"
$ws.Range("C1").Value = "Jun_15"
$ws.Range("D1").Value = "Jun_13"

# Fill the two new data columns (C, D) for every data row with the same "UN"
# placeholder value already used throughout column B.
For ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 3).Value = "UN"
    $ws.Cells.Item($r, 4).Value = "UN"
}
